$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date string in A2: 30.11.2022 -> 30.11.2023
$ws.Range("A2").Value = "30.11.2023"

# Move the active selection from B8 to B5
$ws.Range("B5").Select()
